$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.878.94"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "'3.551.06"
$ws.Range("E3").Value = "  +4.49%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'599.06"
$ws.Range("E5").Value = "  +3.59%  "

$ws.Range("D6").Value = "'136.23"
$ws.Range("E6").Value = "  +3.24%  "

$ws.Range("D7").Value = "'3.549.28"
$ws.Range("E7").Value = "  +4.46%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "'0.496"
$ws.Range("E9").Value = "  +3.64%  "

$ws.Range("E10").Value = "  +3.23%  "

$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("D12").Value = "'0.386"
$ws.Range("E12").Value = "  +4.48%  "

$ws.Range("D13").Value = "'4.157.35"
$ws.Range("E13").Value = "  +4.57%  "

$ws.Range("D14").Value = "'0.0000182"
$ws.Range("E14").Value = "  +3.88%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'3.561.08"
$ws.Range("E15").Value = "  +5.23%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'27.07"
$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "'64.799.27"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").Value = "'10.03"
$ws.Range("E19").Value = "  +7.64%  "

$ws.Range("D20").Value = "'14.40"
$ws.Range("E20").Value = "  +7.85%  "

$ws.Range("D21").Value = "'5.83"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("D22").Value = "'389.41"
$ws.Range("E22").Value = "  +3.46%  "

$ws.Range("D23").Value = "'0.575"
$ws.Range("E23").Value = "  +7.17%  "

$ws.Range("D24").Value = "'3.697.88"
$ws.Range("E24").Value = "  +4.69%  "

$ws.Range("E25").Value = "  +3.94%  "

$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("E27").Value = "  +14.31%  "

$ws.Range("D28").Value = "'7.65"
$ws.Range("E28").Value = "  +8.34%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  +5.63%  "

$ws.Range("E31").Value = "  +5.96%  "

$ws.Range("D32").Value = "'3.561.97"
$ws.Range("E32").Value = "  +4.25%  "

$ws.Range("E33").Value = "  +22.84%  "

$ws.Range("D34").Value = "'23.98"
$ws.Range("E34").Value = "  +5.30%  "

$ws.Range("D36").Value = "'0.145"
$ws.Range("E36").Value = "  +3.44%  "

$ws.Range("D37").Value = "'170.34"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").Value = "'6.93"
$ws.Range("E38").Value = "  +5.49%  "

$ws.Range("D39").Value = "'1.54"
$ws.Range("E39").Value = "  +7.84%  "

$ws.Range("D40").Value = "'4.99"
$ws.Range("E40").Value = "  +9.87%  "

$ws.Range("D41").Value = "'0.0805"
$ws.Range("E41").Value = "  +7.02%  "

$ws.Range("D42").Value = "'0.827"
$ws.Range("E42").Value = "  +4.49%  "

$ws.Range("D43").Value = "'26.77"
$ws.Range("E43").Value = "  +22.48%  "

$ws.Range("D44").Value = "'42.60"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "'4.45"
$ws.Range("E46").Value = "  +5.77%  "

$ws.Range("D47").Value = "'1.20"
$ws.Range("E47").Value = "  +10.08%  "

$ws.Range("E48").Value = "  +4.38%  "

$ws.Range("D49").Value = "'2.452.55"
$ws.Range("E49").Value = "  +12.58%  "

$ws.Range("D50").Value = "'6.88"
$ws.Range("E50").Value = "  +6.99%  "

$ws.Range("E51").Value = "  +17.76%  "
